$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.409481333333334
$ws.Range("H2").Value = 28.228444
$ws.Range("I2").Value = 0.2433300530093958
$ws.Range("J2").Value = 0.2433300530093958
$ws.Range("M2").Value = 31.618405
$ws.Range("N2").Value = 94.855215
$ws.Range("O2").Value = 0.8578613706944929
$ws.Range("P2").Value = 0.8578613706944929
$ws.Range("Q2").Value = 297.5127916372733
$ws.Range("R2").Value = 2677.615124735461
$ws.Range("S2").Value = 0.2087434528058039
$ws.Range("T2").Value = 0.2087434528058039
$ws.Range("G3").Value = 9.409481333333334
$ws.Range("H3").Value = 28.228444
$ws.Range("I3").Value = 0.2433300530093958
$ws.Range("J3").Value = 0.2433300530093958
$ws.Range("O3").Value = 0.08747555172986397
$ws.Range("P3").Value = 0.08747555172986396
$ws.Range("Q3").Value = 30.33718090615689
$ws.Range("R3").Value = 273.034628155412
$ws.Range("S3").Value = 0.02128543063945395
$ws.Range("T3").Value = 0.02128543063945395
$ws.Range("G4").Value = 9.409481333333334
$ws.Range("H4").Value = 28.228444
$ws.Range("I4").Value = 0.2433300530093958
$ws.Range("J4").Value = 0.2433300530093958
$ws.Range("M4").Value = 2.014730333333334
$ws.Range("N4").Value = 6.044191000000001
$ws.Range("O4").Value = 0.05466307757564324
$ws.Range("P4").Value = 0.05466307757564324
$ws.Range("Q4").Value = 18.95756746320045
$ws.Range("R4").Value = 170.618107168804
$ws.Range("S4").Value = 0.01330116956413798
$ws.Range("T4").Value = 0.01330116956413799
$ws.Range("I5").Value = 0.5069354697952918
$ws.Range("J5").Value = 0.5069354697952919
$ws.Range("M5").Value = 31.618405
$ws.Range("N5").Value = 94.855215
$ws.Range("O5").Value = 0.8578613706944929
$ws.Range("P5").Value = 0.8578613706944929
$ws.Range("Q5").Value = 619.8156986096833
$ws.Range("R5").Value = 5578.34128748715
$ws.Range("S5").Value = 0.4348803569722457
$ws.Range("T5").Value = 0.4348803569722458
$ws.Range("I6").Value = 0.5069354697952918
$ws.Range("J6").Value = 0.5069354697952919
$ws.Range("O6").Value = 0.08747555172986397
$ws.Range("P6").Value = 0.08747555172986396
$ws.Range("S6").Value = 0.04434445991178095
$ws.Range("T6").Value = 0.04434445991178095
$ws.Range("I7").Value = 0.5069354697952918
$ws.Range("J7").Value = 0.5069354697952919
$ws.Range("M7").Value = 2.014730333333334
$ws.Range("N7").Value = 6.044191000000001
$ws.Range("O7").Value = 0.05466307757564324
$ws.Range("P7").Value = 0.05466307757564324
$ws.Range("R7").Value = 355.4528889609101
$ws.Range("S7").Value = 0.02771065291126519
$ws.Range("T7").Value = 0.02771065291126519
$ws.Range("G8").Value = 9.657138
$ws.Range("I8").Value = 0.2497344771953123
$ws.Range("J8").Value = 0.2497344771953124
$ws.Range("M8").Value = 31.618405
$ws.Range("N8").Value = 94.855215
$ws.Range("O8").Value = 0.8578613706944929
$ws.Range("P8").Value = 0.8578613706944929
$ws.Range("Q8").Value = 305.34330042489
$ws.Range("R8").Value = 2748.08970382401
$ws.Range("S8").Value = 0.2142375609164432
$ws.Range("T8").Value = 0.2142375609164432
$ws.Range("G9").Value = 9.657138
$ws.Range("I9").Value = 0.2497344771953123
$ws.Range("J9").Value = 0.2497344771953124
$ws.Range("O9").Value = 0.08747555172986397
$ws.Range("P9").Value = 0.08747555172986396
$ws.Range("R9").Value = 280.2208739747219
$ws.Range("S9").Value = 0.02184566117862908
$ws.Range("T9").Value = 0.02184566117862908
$ws.Range("G10").Value = 9.657138
$ws.Range("I10").Value = 0.2497344771953123
$ws.Range("J10").Value = 0.2497344771953124
$ws.Range("M10").Value = 2.014730333333334
$ws.Range("N10").Value = 6.044191000000001
$ws.Range("O10").Value = 0.05466307757564324
$ws.Range("P10").Value = 0.05466307757564324
$ws.Range("S10").Value = 0.01365125510024007
$ws.Range("T10").Value = 0.01365125510024007
